$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The attendance-scanner app re-uploaded this log with the two entries'
# "Log Time" (column D) normalized to a single real Excel time value
# (11:03:15, stored as a numeric time serial with an h:mm:ss display
# format and an explicit black font) instead of the original plain text
# timestamps ("12:54:41" / "12:54:57").
$logTime = (11 * 3600 + 3 * 60 + 15) / 86400.0
$rng = $ws.Range("D2:D3")
$rng.Value = $logTime
$rng.Font.Color = 0
$rng.NumberFormat = "h:mm:ss"

# Matches the selection state ("D2:D3" highlighted) saved in the workbook.
$ws.Range("D2:D3").Select() | Out-Null
